$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 43, shifting existing rows 43-56 down to 44-57
$ws.Rows.Item(43).Insert()

# Fill the newly inserted row 43 with the new record's data
$ws.Range("A43").Value = 4
$ws.Range("B43").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C43").Value = "Los Lagos"
$ws.Range("D43").Value = 45006
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = 10
$ws.Range("F43").Value = 100112012
$ws.Range("G43").Value = "Espinaca"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 25
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 15000
$ws.Range("N43").Value = "`$/cuna 10 kilos"
$ws.Range("O43").Value = "Región Metropolitana"
$ws.Range("P43").Value = 1500
$ws.Range("Q43").Value = 10
$ws.Range("R43").Value = "Hortaliza"
